$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit text format
# applied first, otherwise Excel will silently coerce the string into a Double
# and strip the significant trailing zero / dotted-thousands formatting.
$textCells = @("D5", "D9", "D17", "D18", "D22", "D23", "D25", "D30", "D32", "D35", "D40", "D43", "D45", "D46", "D48", "D49", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.932.36"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.546.35"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "205.77"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").Value = "21.42"
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "1.765.36"
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").Value = "1.546.91"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "26.894.46"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "61.56"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "213.73"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "4.02"
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("D23").Value = "9.16"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  -3.21%  "
$ws.Range("D25").Value = "152.72"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "0.0459"
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("D33").Value = "1.366.17"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("E36").Value = "  +4.97%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  -2.39%  "
$ws.Range("D40").Value = "0.805"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("D43").Value = "5.47"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("D45").Value = "63.39"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "1.75"
$ws.Range("E46").Value = "  -2.65%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.679.31"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("B48").Value = "mCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D48").Value = "2.24"
$ws.Range("E48").Value = "  -3.26%  "
$ws.Range("D49").Value = "86.16"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "0.0950"
$ws.Range("E51").Value = "  -0.09%  "
